$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 16,16
$data[0,0] = [double]"2.231748711887562e-38"
$data[0,1] = [double]"-1"
$data[0,2] = [double]"-7.742867047255508e-33"
$data[0,3] = [double]"-8.778053409629295e-33"
$data[0,4] = [double]"1.000711151416789e-32"
$data[0,5] = [double]"-1.063809148423803e-32"
$data[0,6] = [double]"-5.680688981832213e-12"
$data[0,7] = [double]"5.365095149508155e-12"
$data[0,8] = [double]"5.082721720586695e-12"
$data[0,9] = [double]"4.828585634557375e-12"
$data[0,10] = [double]"2.590589034758215e-22"
$data[0,11] = [double]"-2.260871682934731e-22"
$data[0,12] = [double]"-2.727604027577546e-22"
$data[0,13] = [double]"1.816409443589286e-22"
$data[0,14] = [double]"1.619113835452105e-22"
$data[0,15] = [double]"-1.367827367468036e-22"
$data[1,0] = [double]"-1.676570567569614e-27"
$data[1,1] = [double]"-5.680668107585969e-12"
$data[1,2] = [double]"3.85655271134106e-22"
$data[1,3] = [double]"5.240338119386222e-22"
$data[1,4] = [double]"-7.347910382971129e-22"
$data[1,5] = [double]"4.205674416784645e-25"
$data[1,6] = [double]"1"
$data[1,7] = [double]"1.59373753636074e-15"
$data[1,8] = [double]"3.463733789554758e-16"
$data[1,9] = [double]"-2.685260349080509e-17"
$data[1,10] = [double]"-2.414288686539947e-11"
$data[1,11] = [double]"1.931436229006388e-11"
$data[1,12] = [double]"1.196348392284316e-11"
$data[1,13] = [double]"-1.599932499007658e-11"
$data[1,14] = [double]"-3.735920532742194e-17"
$data[1,15] = [double]"-5.238184556469196e-18"
$data[2,0] = [double]"8.23783031182093e-28"
$data[2,1] = [double]"-5.364996991355424e-12"
$data[2,2] = [double]"2.995545540357545e-22"
$data[2,3] = [double]"3.934106318070613e-22"
$data[2,4] = [double]"4.083130655631688e-24"
$data[2,5] = [double]"7.328767689235939e-22"
$data[2,6] = [double]"2.244109629531476e-15"
$data[2,7] = [double]"-1"
$data[2,8] = [double]"-2.888508131220564e-15"
$data[2,9] = [double]"1.162482032676207e-16"
$data[2,10] = [double]"-3.219086070912584e-11"
$data[2,11] = [double]"-3.8791428784834e-17"
$data[2,12] = [double]"1.292066551498168e-11"
$data[2,13] = [double]"2.10554705060834e-12"
$data[2,14] = [double]"-1.609525049005667e-11"
$data[2,15] = [double]"4.189814855976344e-18"
$data[3,0] = [double]"-3.364924551939469e-20"
$data[3,1] = [double]"-5.080092504873385e-23"
$data[3,2] = [double]"1.073022212612199e-11"
$data[3,3] = [double]"1.207146284773652e-11"
$data[3,4] = [double]"7.922436495644798e-18"
$data[3,5] = [double]"-2.096557156088592e-18"
$data[3,6] = [double]"2.414293677282866e-11"
$data[3,7] = [double]"-3.219053908875856e-11"
$data[3,8] = [double]"-6.935492123487756e-17"
$data[3,9] = [double]"7.768514086828277e-17"
$data[3,10] = [double]"1"
$data[3,11] = [double]"-9.482450287234275e-16"
$data[3,12] = [double]"-5.115920704346041e-16"
$data[3,13] = [double]"1.878791722238532e-16"
$data[3,14] = [double]"7.193873508046173e-16"
$data[3,15] = [double]"5.848981445822339e-17"
$data[4,0] = [double]"-4.99885168993216e-28"
$data[4,1] = [double]"-5.082506165490012e-12"
$data[4,2] = [double]"3.075915050389118e-22"
$data[4,3] = [double]"2.028460970955706e-23"
$data[4,4] = [double]"-4.085231648467161e-22"
$data[4,5] = [double]"5.255744452222151e-22"
$data[4,6] = [double]"8.322573545403726e-16"
$data[4,7] = [double]"2.32184985901952e-15"
$data[4,8] = [double]"-1"
$data[4,9] = [double]"-3.091246208778218e-15"
$data[4,10] = [double]"6.057036818255381e-16"
$data[4,11] = [double]"3.219058922451045e-11"
$data[4,12] = [double]"1.615080698047608e-11"
$data[4,13] = [double]"2.631950403890046e-12"
$data[4,14] = [double]"4.215569406829698e-17"
$data[4,15] = [double]"1.60952930576039e-11"
$data[5,0] = [double]"-1.007155139870415e-22"
$data[5,1] = [double]"-4.724130262788954e-23"
$data[5,2] = [double]"9.657229951783325e-12"
$data[5,3] = [double]"-1.848531944548585e-17"
$data[5,4] = [double]"-1.207143380096243e-11"
$data[5,5] = [double]"-9.549798273697313e-18"
$data[5,6] = [double]"1.931434253825439e-11"
$data[5,7] = [double]"3.410345917924417e-24"
$data[5,8] = [double]"-3.219057089704619e-11"
$data[5,9] = [double]"-2.452461525403291e-24"
$data[5,10] = [double]"-1.259370550597011e-15"
$data[5,11] = [double]"-1"
$data[5,12] = [double]"9.636289446242676e-16"
$data[5,13] = [double]"-5.534092425119255e-16"
$data[5,14] = [double]"1.053792632085691e-15"
$data[5,15] = [double]"3.391134320120824e-16"
$data[6,0] = [double]"3.093488856205732e-17"
$data[6,1] = [double]"-4.386088127070814e-23"
$data[6,2] = [double]"8.779164453220413e-12"
$data[6,3] = [double]"7.327766258857363e-17"
$data[6,4] = [double]"1.432019698021412e-16"
$data[6,5] = [double]"1.207131452158868e-11"
$data[6,6] = [double]"8.826464305920726e-26"
$data[6,7] = [double]"-1.931434253822927e-11"
$data[6,8] = [double]"-2.414292817278689e-11"
$data[6,9] = [double]"1.034453723025051e-25"
$data[6,10] = [double]"-1.109909007560883e-15"
$data[6,11] = [double]"-1.516375350366841e-15"
$data[6,12] = [double]"-0.6689647316216069"
$data[6,13] = [double]"-0.1090164664112791"
$data[6,14] = [double]"-5.0146187834951e-17"
$data[6,15] = [double]"5.800604126749346e-16"
$data[7,0] = [double]"4.389624060269413e-12"
$data[7,1] = [double]"-1.141901946943713e-33"
$data[7,2] = [double]"1"
$data[7,3] = [double]"-1.397033491803106e-15"
$data[7,4] = [double]"1.070798920892243e-16"
$data[7,5] = [double]"-1.452738233480887e-17"
$data[7,6] = [double]"-8.39351447153323e-22"
$data[7,7] = [double]"8.289799725241271e-22"
$data[7,8] = [double]"7.771777564969107e-22"
$data[7,9] = [double]"3.226394493872073e-27"
$data[7,10] = [double]"-1.073024865693331e-11"
$data[7,11] = [double]"9.657193530485186e-12"
$data[7,12] = [double]"5.873286546139251e-12"
$data[7,13] = [double]"9.569053443397464e-13"
$data[7,14] = [double]"3.478292783793902e-18"
$data[7,15] = [double]"6.138332075495001e-18"
$data[8,0] = [double]"-2.709834350789274e-27"
$data[8,1] = [double]"-4.828325407330431e-12"
$data[8,2] = [double]"-3.067266319856915e-23"
$data[8,3] = [double]"2.439254143947731e-22"
$data[8,4] = [double]"-3.114995750692445e-22"
$data[8,5] = [double]"3.934873833296838e-22"
$data[8,6] = [double]"3.689603067076522e-16"
$data[8,7] = [double]"-2.482392588656102e-16"
$data[8,8] = [double]"2.638685493388011e-15"
$data[8,9] = [double]"-1"
$data[8,10] = [double]"-5.072008259540997e-16"
$data[8,11] = [double]"7.121684356580467e-19"
$data[8,12] = [double]"2.39270193881347e-11"
$data[8,13] = [double]"-3.199868430412053e-11"
$data[8,14] = [double]"-2.414296792954197e-11"
$data[8,15] = [double]"1.931433613125746e-11"
$data[9,0] = [double]"3.990565576023635e-17"
$data[9,1] = [double]"-4.411731537076132e-23"
$data[9,2] = [double]"-6.275049293301945e-17"
$data[9,3] = [double]"9.657073033623012e-12"
$data[9,4] = [double]"-1.073019920959177e-11"
$data[9,5] = [double]"-7.515099769079732e-17"
$data[9,6] = [double]"1.609528544852462e-11"
$data[9,7] = [double]"-1.030582811001717e-26"
$data[9,8] = [double]"5.666405316807146e-25"
$data[9,9] = [double]"-3.219057089704946e-11"
$data[9,10] = [double]"1.674625396984456e-16"
$data[9,11] = [double]"-2.761973912883562e-17"
$data[9,12] = [double]"-0.7432941462479248"
$data[9,13] = [double]"0.9940399438911891"
$data[9,14] = [double]"-6.405068774807477e-16"
$data[9,15] = [double]"5.814498672251493e-16"
$data[10,0] = [double]"2.318992847370914e-20"
$data[10,1] = [double]"-4.101484788679566e-23"
$data[10,2] = [double]"-6.26126423518508e-17"
$data[10,3] = [double]"8.779264969891236e-12"
$data[10,4] = [double]"-2.996622578430979e-17"
$data[10,5] = [double]"1.073019912382238e-11"
$data[10,6] = [double]"-4.300665460244105e-22"
$data[10,7] = [double]"-1.609528544858974e-11"
$data[10,8] = [double]"1.396007887053938e-23"
$data[10,9] = [double]"-2.414292817279615e-11"
$data[10,10] = [double]"-6.810427651453845e-16"
$data[10,11] = [double]"8.867588615505553e-16"
$data[10,12] = [double]"-3.26704854165604e-16"
$data[10,13] = [double]"5.905304717420313e-16"
$data[10,14] = [double]"1"
$data[10,15] = [double]"2.012757527758671e-15"
$data[11,0] = [double]"4.198762716307365e-12"
$data[11,1] = [double]"-1.009985749897254e-33"
$data[11,2] = [double]"9.636416885034894e-16"
$data[11,3] = [double]"1"
$data[11,4] = [double]"-5.36500094496427e-17"
$data[11,5] = [double]"4.458093288320195e-16"
$data[11,6] = [double]"-9.714667506989633e-22"
$data[11,7] = [double]"9.326104496450015e-22"
$data[11,8] = [double]"2.141328404866438e-27"
$data[11,9] = [double]"7.771726278103482e-22"
$data[11,10] = [double]"-1.207146409970206e-11"
$data[11,11] = [double]"-6.913053972440056e-20"
$data[11,12] = [double]"7.178117670797212e-12"
$data[11,13] = [double]"-9.599540615954022e-12"
$data[11,14] = [double]"-8.779246611389223e-12"
$data[11,15] = [double]"2.48604860623767e-23"
$data[12,0] = [double]"-1.363247918779208e-20"
$data[12,1] = [double]"-3.842194312588082e-23"
$data[12,2] = [double]"1.96684543438495e-17"
$data[12,3] = [double]"-6.209438931723159e-18"
$data[12,4] = [double]"-8.779236340784484e-12"
$data[12,5] = [double]"9.657168217908482e-12"
$data[12,6] = [double]"-5.796898387193339e-20"
$data[12,7] = [double]"-2.219195522183423e-20"
$data[12,8] = [double]"-1.609528543151388e-11"
$data[12,9] = [double]"-1.931434255719504e-11"
$data[12,10] = [double]"3.048202111996052e-17"
$data[12,11] = [double]"-2.127441337955986e-16"
$data[12,12] = [double]"-8.65273874625468e-16"
$data[12,13] = [double]"2.535015456704221e-16"
$data[12,14] = [double]"3.232611405306137e-15"
$data[12,15] = [double]"-1"
$data[13,0] = [double]"4.023897749873957e-12"
$data[13,1] = [double]"-8.935451958260775e-34"
$data[13,2] = [double]"2.021477049285494e-16"
$data[13,3] = [double]"-1.048249412700043e-16"
$data[13,4] = [double]"-1"
$data[13,5] = [double]"-1.570680296142588e-15"
$data[13,6] = [double]"-1.13985613885977e-21"
$data[13,7] = [double]"3.170345555754564e-31"
$data[13,8] = [double]"9.326095690807819e-22"
$data[13,9] = [double]"8.289862837575212e-22"
$data[13,10] = [double]"8.882880721229753e-22"
$data[13,11] = [double]"1.207146409776064e-11"
$data[13,12] = [double]"7.975790826044586e-12"
$data[13,13] = [double]"-1.066638811050178e-11"
$data[13,14] = [double]"3.587011519156772e-21"
$data[13,15] = [double]"8.779246608354961e-12"
$data[14,0] = [double]"3.862625908880551e-12"
$data[14,1] = [double]"-7.920831019283158e-34"
$data[14,2] = [double]"3.194773389598348e-17"
$data[14,3] = [double]"-2.381594983173658e-16"
$data[14,4] = [double]"-9.822805366194853e-16"
$data[14,5] = [double]"1"
$data[14,6] = [double]"9.781508744603279e-31"
$data[14,7] = [double]"1.139856140506642e-21"
$data[14,8] = [double]"9.714683009923193e-22"
$data[14,9] = [double]"8.393486126241836e-22"
$data[14,10] = [double]"-8.627754183637854e-22"
$data[14,11] = [double]"-3.955176453598527e-21"
$data[14,12] = [double]"8.075010142429449e-12"
$data[14,13] = [double]"1.316200083213747e-12"
$data[14,14] = [double]"-1.073019029912856e-11"
$data[14,15] = [double]"9.657171269116426e-12"
$data[15,0] = [double]"1"
$data[15,1] = [double]"3.707189766630389e-44"
$data[15,2] = [double]"-4.389623327864951e-12"
$data[15,3] = [double]"-4.198770128267268e-12"
$data[15,4] = [double]"4.023819818031322e-12"
$data[15,5] = [double]"-3.862875931147207e-12"
$data[15,6] = [double]"1.055404632371026e-32"
$data[15,7] = [double]"-1.000542051913699e-32"
$data[15,8] = [double]"-8.93452432511236e-33"
$data[15,9] = [double]"-7.856347481115003e-33"
$data[15,10] = [double]"7.028571994197634e-23"
$data[15,11] = [double]"-6.903515939756928e-23"
$data[15,12] = [double]"6.394251278974946e-17"
$data[15,13] = [double]"-3.226509235403771e-17"
$data[15,14] = [double]"5.510613654697584e-23"
$data[15,15] = [double]"-5.087355226303406e-23"

$ws.Range("A2:P17").Value = $data
